$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $ws.Columns.Item(4).Width = 100
  Write-Host "Width set succeeded"
} catch {
  Write-Host "Width set failed: $_"
}
